$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.57077
$ws.Range("H2").Value = 4.71231
$ws.Range("I2").Value = 0.02582502173444737
$ws.Range("J2").Value = 0.02582502173444737
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.177232
$ws.Range("N2").Value = 0.531696
$ws.Range("O2").Value = 0.0005104719838156216
$ws.Range("P2").Value = 0.0005104719838156217
$ws.Range("Q2").Value = 0.27839070864
$ws.Range("R2").Value = 2.50551637776
$ws.Range("S2").Value = 0.00001318295007686489
$ws.Range("T2").Value = 0.0000131829500768649

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.57077
$ws.Range("H3").Value = 4.71231
$ws.Range("I3").Value = 0.02582502173444737
$ws.Range("J3").Value = 0.02582502173444737
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.08113566666666668
$ws.Range("N3").Value = 0.243407
$ws.Range("O3").Value = 0.0002336907822601807
$ws.Range("P3").Value = 0.0002336907822601807
$ws.Range("Q3").Value = 0.12744547113
$ws.Range("R3").Value = 1.14700924017
$ws.Range("S3").Value = 0.000006035069531009176
$ws.Range("T3").Value = 0.000006035069531009176

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.57077
$ws.Range("H4").Value = 4.71231
$ws.Range("I4").Value = 0.02582502173444737
$ws.Range("J4").Value = 0.02582502173444737
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 274.5137023333334
$ws.Range("N4").Value = 823.541107
$ws.Range("O4").Value = 0.7906673411949746
$ws.Range("P4").Value = 0.7906673411949746
$ws.Range("Q4").Value = 431.1978882141301
$ws.Range("R4").Value = 3880.780993927171
$ws.Range("S4").Value = 0.02041900127107793
$ws.Range("T4").Value = 0.02041900127107793

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.57077
$ws.Range("H5").Value = 4.71231
$ws.Range("I5").Value = 0.02582502173444737
$ws.Range("J5").Value = 0.02582502173444737
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.042388
$ws.Range("N5").Value = 0.127164
$ws.Range("O5").Value = 0.0001220879211991998
$ws.Range("P5").Value = 0.0001220879211991998
$ws.Range("Q5").Value = 0.06658179876000002
$ws.Range("R5").Value = 0.5992361888400001
$ws.Range("S5").Value = 0.000003152923218482832
$ws.Range("T5").Value = 0.000003152923218482832

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.57077
$ws.Range("H6").Value = 4.71231
$ws.Range("I6").Value = 0.02582502173444737
$ws.Range("J6").Value = 0.02582502173444737
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 72.37795533333333
$ws.Range("N6").Value = 217.133866
$ws.Range("O6").Value = 0.2084664081177503
$ws.Range("P6").Value = 0.2084664081177503
$ws.Range("Q6").Value = 113.68912089894
$ws.Range("R6").Value = 1023.20208809046
$ws.Range("S6").Value = 0.005383649520543078
$ws.Range("T6").Value = 0.005383649520543078

# Row 7
$ws.Range("I7").Value = 0.934831682683009
$ws.Range("J7").Value = 0.934831682683009
$ws.Range("M7").Value = 0.177232
$ws.Range("N7").Value = 0.531696
$ws.Range("O7").Value = 0.0005104719838156216
$ws.Range("P7").Value = 0.0005104719838156217
$ws.Range("Q7").Value = 10.07737601452267
$ws.Range("R7").Value = 90.696384130704
$ws.Range("S7").Value = 0.0004772053835928913
$ws.Range("T7").Value = 0.0004772053835928914

# Row 8
$ws.Range("I8").Value = 0.934831682683009
$ws.Range("J8").Value = 0.934831682683009
$ws.Range("O8").Value = 0.0002336907822601807
$ws.Range("P8").Value = 0.0002336907822601807
$ws.Range("S8").Value = 0.0002184615472077934
$ws.Range("T8").Value = 0.0002184615472077934

# Row 9
$ws.Range("I9").Value = 0.934831682683009
$ws.Range("J9").Value = 0.934831682683009
$ws.Range("M9").Value = 274.5137023333334
$ws.Range("N9").Value = 823.541107
$ws.Range("O9").Value = 0.7906673411949746
$ws.Range("P9").Value = 0.7906673411949746
$ws.Range("Q9").Value = 15608.7941204283
$ws.Range("R9").Value = 140479.1470838547
$ws.Range("S9").Value = 0.7391408810117989
$ws.Range("T9").Value = 0.7391408810117989

# Row 10
$ws.Range("I10").Value = 0.934831682683009
$ws.Range("J10").Value = 0.934831682683009
$ws.Range("M10").Value = 0.042388
$ws.Range("N10").Value = 0.127164
$ws.Range("O10").Value = 0.0001220879211991998
$ws.Range("P10").Value = 0.0001220879211991998
$ws.Range("Q10").Value = 2.410173188270667
$ws.Range("R10").Value = 21.691558694436
$ws.Range("S10").Value = 0.0001141316568099185
$ws.Range("T10").Value = 0.0001141316568099185

# Row 11
$ws.Range("I11").Value = 0.934831682683009
$ws.Range("J11").Value = 0.934831682683009
$ws.Range("M11").Value = 72.37795533333333
$ws.Range("N11").Value = 217.133866
$ws.Range("O11").Value = 0.2084664081177503
$ws.Range("P11").Value = 0.2084664081177503
$ws.Range("Q11").Value = 4115.396040536282
$ws.Range("R11").Value = 37038.56436482654
$ws.Range("S11").Value = 0.1948810030835994
$ws.Range("T11").Value = 0.1948810030835994

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.9273763333333335
$ws.Range("H12").Value = 2.782129
$ws.Range("I12").Value = 0.01524698967025436
$ws.Range("J12").Value = 0.01524698967025436
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.177232
$ws.Range("N12").Value = 0.531696
$ws.Range("O12").Value = 0.0005104719838156216
$ws.Range("P12").Value = 0.0005104719838156217
$ws.Range("Q12").Value = 0.1643607623093333
$ws.Range("R12").Value = 1.479246860784
$ws.Range("S12").Value = 0.000007783161064191034
$ws.Range("T12").Value = 0.000007783161064191035

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.9273763333333335
$ws.Range("H13").Value = 2.782129
$ws.Range("I13").Value = 0.01524698967025436
$ws.Range("J13").Value = 0.01524698967025436
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.08113566666666668
$ws.Range("N13").Value = 0.243407
$ws.Range("O13").Value = 0.0002336907822601807
$ws.Range("P13").Value = 0.0002336907822601807
$ws.Range("Q13").Value = 0.0752432970558889
$ws.Range("R13").Value = 0.6771896735030001
$ws.Range("S13").Value = 0.000003563080943154637
$ws.Range("T13").Value = 0.000003563080943154637

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.9273763333333335
$ws.Range("H14").Value = 2.782129
$ws.Range("I14").Value = 0.01524698967025436
$ws.Range("J14").Value = 0.01524698967025436
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 274.5137023333334
$ws.Range("N14").Value = 823.541107
$ws.Range("O14").Value = 0.7906673411949746
$ws.Range("P14").Value = 0.7906673411949746
$ws.Range("Q14").Value = 254.5775107196448
$ws.Range("R14").Value = 2291.197596476803
$ws.Range("S14").Value = 0.01205529678380726
$ws.Range("T14").Value = 0.01205529678380726

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.9273763333333335
$ws.Range("H15").Value = 2.782129
$ws.Range("I15").Value = 0.01524698967025436
$ws.Range("J15").Value = 0.01524698967025436
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 0.6666666666666666
$ws.Range("M15").Value = 0.042388
$ws.Range("N15").Value = 0.127164
$ws.Range("O15").Value = 0.0001220879211991998
$ws.Range("P15").Value = 0.0001220879211991998
$ws.Range("Q15").Value = 0.03930962801733334
$ws.Range("R15").Value = 0.353786652156
$ws.Range("S15").Value = 0.000001861473273387027
$ws.Range("T15").Value = 0.000001861473273387027

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.9273763333333335
$ws.Range("H16").Value = 2.782129
$ws.Range("I16").Value = 0.01524698967025436
$ws.Range("J16").Value = 0.01524698967025436
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 72.37795533333333
$ws.Range("N16").Value = 217.133866
$ws.Range("O16").Value = 0.2084664081177503
$ws.Range("P16").Value = 0.2084664081177503
$ws.Range("Q16").Value = 67.12160283119046
$ws.Range("R16").Value = 604.0944254807141
$ws.Range("S16").Value = 0.003178485171166369
$ws.Range("T16").Value = 0.003178485171166369

# Row 17
$ws.Range("G17").Value = 0.7810079999999999
$ws.Range("H17").Value = 2.343024
$ws.Range("I17").Value = 0.0128405486320577
$ws.Range("J17").Value = 0.0128405486320577
$ws.Range("M17").Value = 0.177232
$ws.Range("N17").Value = 0.531696
$ws.Range("O17").Value = 0.0005104719838156216
$ws.Range("P17").Value = 0.0005104719838156217
$ws.Range("Q17").Value = 0.138419609856
$ws.Range("R17").Value = 1.245776488704
$ws.Range("S17").Value = 0.000006554740333487459
$ws.Range("T17").Value = 0.00000655474033348746

# Row 18
$ws.Range("G18").Value = 0.7810079999999999
$ws.Range("H18").Value = 2.343024
$ws.Range("I18").Value = 0.0128405486320577
$ws.Range("J18").Value = 0.0128405486320577
$ws.Range("O18").Value = 0.0002336907822601807
$ws.Range("P18").Value = 0.0002336907822601807
$ws.Range("Q18").Value = 0.063367604752
$ws.Range("R18").Value = 0.570308442768
$ws.Range("S18").Value = 0.000003000717854475457
$ws.Range("T18").Value = 0.000003000717854475457

# Row 19
$ws.Range("G19").Value = 0.7810079999999999
$ws.Range("H19").Value = 2.343024
$ws.Range("I19").Value = 0.0128405486320577
$ws.Range("J19").Value = 0.0128405486320577
$ws.Range("M19").Value = 274.5137023333334
$ws.Range("N19").Value = 823.541107
$ws.Range("O19").Value = 0.7906673411949746
$ws.Range("P19").Value = 0.7906673411949746
$ws.Range("Q19").Value = 214.397397631952
$ws.Range("R19").Value = 1929.576578687568
$ws.Range("S19").Value = 0.01015260244639383
$ws.Range("T19").Value = 0.01015260244639383

# Row 20
$ws.Range("G20").Value = 0.7810079999999999
$ws.Range("H20").Value = 2.343024
$ws.Range("I20").Value = 0.0128405486320577
$ws.Range("J20").Value = 0.0128405486320577
$ws.Range("M20").Value = 0.042388
$ws.Range("N20").Value = 0.127164
$ws.Range("O20").Value = 0.0001220879211991998
$ws.Range("P20").Value = 0.0001220879211991998
$ws.Range("Q20").Value = 0.033105367104
$ws.Range("R20").Value = 0.297948303936
$ws.Range("S20").Value = 0.000001567675889545152
$ws.Range("T20").Value = 0.000001567675889545152

# Row 21
$ws.Range("G21").Value = 0.7810079999999999
$ws.Range("H21").Value = 2.343024
$ws.Range("I21").Value = 0.0128405486320577
$ws.Range("J21").Value = 0.0128405486320577
$ws.Range("M21").Value = 72.37795533333333
$ws.Range("N21").Value = 217.133866
$ws.Range("O21").Value = 0.2084664081177503
$ws.Range("P21").Value = 0.2084664081177503
$ws.Range("Q21").Value = 56.52776213897599
$ws.Range("R21").Value = 508.749859250784
$ws.Range("S21").Value = 0.00267682305158636
$ws.Range("T21").Value = 0.00267682305158636

# Row 22
$ws.Range("G22").Value = 0.6846153333333334
$ws.Range("H22").Value = 2.053846
$ws.Range("I22").Value = 0.01125575728023152
$ws.Range("J22").Value = 0.01125575728023152
$ws.Range("M22").Value = 0.177232
$ws.Range("N22").Value = 0.531696
$ws.Range("O22").Value = 0.0005104719838156216
$ws.Range("P22").Value = 0.0005104719838156217
$ws.Range("Q22").Value = 0.1213357447573333
$ws.Range("R22").Value = 1.092021702816
$ws.Range("S22").Value = 0.000005745748748186909
$ws.Range("T22").Value = 0.00000574574874818691

# Row 23
$ws.Range("G23").Value = 0.6846153333333334
$ws.Range("H23").Value = 2.053846
$ws.Range("I23").Value = 0.01125575728023152
$ws.Range("J23").Value = 0.01125575728023152
$ws.Range("O23").Value = 0.0002336907822601807
$ws.Range("P23").Value = 0.0002336907822601807
$ws.Range("Q23").Value = 0.05554672148022223
$ws.Range("R23").Value = 0.4999204933220001
$ws.Range("S23").Value = 0.000002630366723748028
$ws.Range("T23").Value = 0.000002630366723748028

# Row 24
$ws.Range("G24").Value = 0.6846153333333334
$ws.Range("H24").Value = 2.053846
$ws.Range("I24").Value = 0.01125575728023152
$ws.Range("J24").Value = 0.01125575728023152
$ws.Range("M24").Value = 274.5137023333334
$ws.Range("N24").Value = 823.541107
$ws.Range("O24").Value = 0.7906673411949746
$ws.Range("P24").Value = 0.7906673411949746
$ws.Range("Q24").Value = 187.9362898275025
$ws.Range("R24").Value = 1691.426608447522
$ws.Range("S24").Value = 0.008899559681896634
$ws.Range("T24").Value = 0.008899559681896634

# Row 25
$ws.Range("G25").Value = 0.6846153333333334
$ws.Range("H25").Value = 2.053846
$ws.Range("I25").Value = 0.01125575728023152
$ws.Range("J25").Value = 0.01125575728023152
$ws.Range("M25").Value = 0.042388
$ws.Range("N25").Value = 0.127164
$ws.Range("O25").Value = 0.0001220879211991998
$ws.Range("P25").Value = 0.0001220879211991998
$ws.Range("Q25").Value = 0.02901947474933333
$ws.Range("R25").Value = 0.261175272744
$ws.Range("S25").Value = 0.000001374192007866225
$ws.Range("T25").Value = 0.000001374192007866225

# Row 26
$ws.Range("G26").Value = 0.6846153333333334
$ws.Range("H26").Value = 2.053846
$ws.Range("I26").Value = 0.01125575728023152
$ws.Range("J26").Value = 0.01125575728023152
$ws.Range("M26").Value = 72.37795533333333
$ws.Range("N26").Value = 217.133866
$ws.Range("O26").Value = 0.2084664081177503
$ws.Range("P26").Value = 0.2084664081177503
$ws.Range("Q26").Value = 49.55105801651511
$ws.Range("R26").Value = 445.959522148636
$ws.Range("S26").Value = 0.002346447290855083
$ws.Range("T26").Value = 0.002346447290855083

